$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("ThuongHieu") entirely, shifting C->B, D->C
$ws.Range("B:B").Delete()

# Update the selection to match the target state
$ws.Range("D9").Select()
